# Update gh-pages to output generated at 456a3b4
# Applies updated "want-to-go" counts (F column) and lowest price (G5)
# to both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of cell address -> new value, applied identically to both sheets
$updates = @{
    "F2"  = 156
    "F3"  = 1778
    "G5"  = 68
    "F11" = 21
    "F12" = 84
    "F19" = 216
    "F20" = 34
    "F21" = 436
    "F22" = 335
    "F23" = 129
    "F27" = 721
    "F28" = 2520
    "F31" = 504
    "F32" = 820
    "F36" = 375
    "F38" = 574
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
